# Update "想去人数" (F column) figures across the relevant worksheets
# to reflect the newly generated output (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 89
$wsExhibit.Range("F3").Value = 11930
$wsExhibit.Range("F4").Value = 25
$wsExhibit.Range("F5").Value = 221
$wsExhibit.Range("F6").Value = 357
$wsExhibit.Range("F7").Value = 226
$wsExhibit.Range("F8").Value = 11829
$wsExhibit.Range("F9").Value = 492
$wsExhibit.Range("F10").Value = 1175
$wsExhibit.Range("F11").Value = 104
$wsExhibit.Range("F12").Value = 67
$wsExhibit.Range("F13").Value = 1783
$wsExhibit.Range("F14").Value = 5862
$wsExhibit.Range("F15").Value = 126
$wsExhibit.Range("F16").Value = 3541

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 5

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 89
$wsAll.Range("F5").Value = 11930
$wsAll.Range("F6").Value = 25
$wsAll.Range("F7").Value = 221
$wsAll.Range("F8").Value = 5
$wsAll.Range("F9").Value = 357
$wsAll.Range("F10").Value = 226
$wsAll.Range("F11").Value = 11829
$wsAll.Range("F12").Value = 492
$wsAll.Range("F13").Value = 1175
$wsAll.Range("F14").Value = 104
$wsAll.Range("F15").Value = 67
$wsAll.Range("F16").Value = 1783
$wsAll.Range("F18").Value = 5862
$wsAll.Range("F19").Value = 126
$wsAll.Range("F20").Value = 3541
